# Update column G ("K") values for rows 2-42 on the active sheet.
# This corresponds to regenerating save_data to use K (strikeouts) instead
# of Strike# in the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 3
    22 = 1
    23 = 1
    24 = 2
    25 = 4
    26 = 1
    27 = 1
    28 = 3
    29 = 1
    30 = 1
    31 = 2
    32 = 0
    33 = 1
    34 = 1
    35 = 3
    36 = 1
    37 = 4
    38 = 2
    39 = 1
    40 = 2
    41 = 1
    42 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
